$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.763.09"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "1.889.86"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'0.7796"
$ws.Range("E5").Value = "  -5.85%  "
$ws.Range("D6").Value = "'241.80"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("D8").Value = "'0.3163"
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "'25.32"
$ws.Range("E9").Value = "  -4.69%  "
$ws.Range("D10").Value = "'0.07019"
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").Value = "'0.08038"
$ws.Range("D12").Value = "'0.7656"
$ws.Range("E12").Value = "  +2.59%  "
$ws.Range("D13").Value = "1.905.37"
$ws.Range("E13").Value = "  +0.16%  "
$ws.Range("D14").Value = "'5.278"
$ws.Range("E14").Value = "  +1.60%  "
$ws.Range("D15").Value = "'91.92"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").Value = "29.789.18"
$ws.Range("E16").Value = "  -0.63%  "
$ws.Range("D17").Value = "'13.82"
$ws.Range("E17").Value = "  -1.43%  "
$ws.Range("D18").Value = "'5.897"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("D19").Value = "'242.47"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").Value = "'0.000007701"
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").Value = "'8.232"
$ws.Range("E21").Value = "  +18.75%  "
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "2.149.71"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'0.1647"
$ws.Range("E25").Value = "  +4.24%  "
$ws.Range("D26").Value = "'9.291"
$ws.Range("E26").Value = "  +1.41%  "
$ws.Range("D27").Value = "'165.20"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "'18.65"
$ws.Range("E28").Value = "  -0.96%  "
$ws.Range("D29").Value = "'2.045"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").Value = "'1.401"
$ws.Range("E30").Value = "  +1.65%  "
$ws.Range("D31").Value = "'1.536"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("D32").Value = "'4.395"
$ws.Range("E32").Value = "  +3.45%  "
$ws.Range("D33").Value = "'0.05610"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "'4.030"
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("D35").Value = "'1.259"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "'0.7344"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").Value = "'1.005"
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D39").Value = "'0.01904"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").Value = "'2.763"
$ws.Range("E40").Value = "  -0.74%  "
$ws.Range("D41").Value = "'0.4395"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").Value = "'72.33"
$ws.Range("E42").Value = "  +0.83%  "
$ws.Range("D43").Value = "'5.803"
$ws.Range("E43").Value = "  -2.52%  "
$ws.Range("D44").Value = "'0.9997"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "'0.8377"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").Value = "'102.15"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "1.017.05"
$ws.Range("E47").Value = "  +3.17%  "
$ws.Range("D48").Value = "'1.860"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "'9.884"
$ws.Range("E49").Value = "  +2.02%  "
$ws.Range("D50").Value = "'7.395"
$ws.Range("E50").Value = "  -2.40%  "
$ws.Range("D51").Value = "2.054.39"
$ws.Range("E51").Value = "  -0.07%  "
